$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - logistic_embeddings
$ws.Range("C5").Value = 0.315
$ws.Range("D5").Value = 0.442
$ws.Range("E5").Value = 0.48
$ws.Range("F5").Value = 0.542
$ws.Range("G5").Value = 0.5580000000000001
$ws.Range("H5").Value = 0.576

# Row 7 - classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.315
$ws.Range("E7").Value = 0.48
$ws.Range("F7").Value = 0.542

# Row 8 - BERT-base
$ws.Range("C8").Value = 0.309
$ws.Range("D8").Value = 0.524
$ws.Range("E8").Value = 0.5649999999999999
$ws.Range("F8").Value = 0.608
$ws.Range("G8").Value = 0.641
$ws.Range("H8").Value = 0.655

# Row 9 - BERT-base-nli
$ws.Range("B9").Value = 0.291
$ws.Range("C9").Value = 0.419
$ws.Range("D9").Value = 0.549
$ws.Range("E9").Value = 0.577
$ws.Range("F9").Value = 0.593
$ws.Range("G9").Value = 0.625
$ws.Range("H9").Value = 0.638
